# Mark up additional submissions as "Yes" (marked) for rows 14-19.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G14").Value = "Yes"

$ws.Range("G15").Value = "Yes"
$ws.Range("G16").Value = "Yes"
$ws.Range("G17").Value = "Yes"
$ws.Range("G18").Value = "Yes"
$ws.Range("G19").Value = "Yes"

# Eduard Zaharia (row 15) gets a mailto hyperlink on his email cell, like
# the other "emailed directly" entries (e.g. row 20).
$ws.Hyperlinks.Add($ws.Range("B15"), "mailto:ezaharia1@sheffield.ac.uk")
$ws.Range("B15").Style = "Hyperlink"

# A few cells got re-clicked/re-linked and picked up the standard
# "Hyperlink" cell style along the way.
$ws.Range("E15").Style = "Hyperlink"
$ws.Range("E17").Style = "Hyperlink"
$ws.Range("E19").Style = "Hyperlink"

# Notes on marking progress.
$ws.Range("G43").Value = "28 can mark"
$ws.Range("G45").Value = "29 done by Tues (2 per day)"
$ws.Range("G44").Value = "21 done by today"
$ws.Range("G47").Value = "10 remaining over 3 days"
$ws.Range("H47").Value = "(could be more drawn out than that)"

$ws.Range("G15").Select() | Out-Null
